# Realestate Update resale numbers 2025-01-25 12:01
# Appends a new data row (row 42) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

# Columns A-D hold text (date/time/weekday/week-code) values in the source
# data. Force a Text number format first so Excel's COM layer doesn't
# auto-coerce date-looking / zero-padded strings into dates or numbers.
$ws.Range("A" + $row + ":D" + $row).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-25"   # A: Date
$ws.Cells.Item($row, 2).Value = "12:01:25"     # B: Time
$ws.Cells.Item($row, 3).Value = "Saturday"     # C: Weekday
$ws.Cells.Item($row, 4).Value = "03"           # D: Week

# Columns E-T hold numeric resale figures (-1 denotes missing data).
$ws.Cells.Item($row, 5).Value  = 126167   # E: Beijing
$ws.Cells.Item($row, 6).Value  = 142081   # F: Guangzhou
$ws.Cells.Item($row, 7).Value  = 168276   # G: Suzhou
$ws.Cells.Item($row, 8).Value  = 158544   # H: Hangzhou
$ws.Cells.Item($row, 9).Value  = -1       # I: Nanjing
$ws.Cells.Item($row, 10).Value = 142690   # J: Xi_an
$ws.Cells.Item($row, 11).Value = -1       # K: Chengdu
$ws.Cells.Item($row, 12).Value = -1       # L: Chongqing
$ws.Cells.Item($row, 13).Value = 191589   # M: Tianjin
$ws.Cells.Item($row, 14).Value = 115644   # N: Hefei
$ws.Cells.Item($row, 15).Value = 45529    # O: Fuzhou
$ws.Cells.Item($row, 16).Value = 28400    # P: Xiamen
$ws.Cells.Item($row, 17).Value = 65080    # Q: Changsha
$ws.Cells.Item($row, 18).Value = -1       # R: Shanghai
$ws.Cells.Item($row, 19).Value = 46380    # S: Shenzhen
$ws.Cells.Item($row, 20).Value = -1       # T: Wuhan
